$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing data-row formatting (centered, like rows 2-32) down to the new rows 33-37
# so the newly appended rows match the look of the rest of the table.
$ws.Range("A2:C2").Copy()
$ws.Range("A33:C37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rewrite the Variable Name / Type / Notes table ---
# Row 29 ("star_scores") is dropped, "Merged Budget"/"Merged Revenue" shift up one row,
# and the old "director_scores" row is replaced by the new star/production-company
# scoring breakdown rows (31-37).
$ws.Range("A1").Value = "Variable Name"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Notes"
$ws.Range("A2").Value = "Unnamed: 0"
$ws.Range("B2").Value = "Categorical"
$ws.Range("C2").Value = "Index"
$ws.Range("A3").Value = "Running Time"
$ws.Range("B3").Value = "Numerical"
$ws.Range("C3").Value = "hrs mins"
$ws.Range("A4").Value = "Genres"
$ws.Range("B4").Value = "Categorical"
$ws.Range("A5").Value = "IMDB Title"
$ws.Range("B5").Value = "Categorical"
$ws.Range("A6").Value = "MPAA"
$ws.Range("B6").Value = "Categorical"
$ws.Range("C6").Value = "Maturity Rating"
$ws.Range("A7").Value = "Domestic Distributor"
$ws.Range("B7").Value = "Categorical"
$ws.Range("A8").Value = "Domestic Opening"
$ws.Range("B8").Value = "Categorical"
$ws.Range("A9").Value = "Earliest Release Date"
$ws.Range("B9").Value = "Numerical"
$ws.Range("A10").Value = "IMDB Budget"
$ws.Range("B10").Value = "Numerical"
$ws.Range("A11").Value = "TMDB Budget"
$ws.Range("B11").Value = "Numerical"
$ws.Range("A12").Value = "Genre ID"
$ws.Range("B12").Value = "Categorical"
$ws.Range("A13").Value = "Genre Name"
$ws.Range("B13").Value = "Categorical"
$ws.Range("A14").Value = "IMDB ID"
$ws.Range("B14").Value = "Categorical"
$ws.Range("A15").Value = "Production Company ID"
$ws.Range("B15").Value = "Categorical"
$ws.Range("A16").Value = "Production Company Name"
$ws.Range("B16").Value = "Categorical"
$ws.Range("A17").Value = "Release Date"
$ws.Range("B17").Value = "Numerical"
$ws.Range("A18").Value = "TMDB Revenue"
$ws.Range("B18").Value = "Numerical"
$ws.Range("A19").Value = "Runtime"
$ws.Range("B19").Value = "Numerical"
$ws.Range("C19").Value = "Integer"
$ws.Range("A20").Value = "TMDB Title"
$ws.Range("B20").Value = "Categorical"
$ws.Range("A21").Value = "TMDB ID"
$ws.Range("B21").Value = "Categorical"
$ws.Range("A22").Value = "cast_ids"
$ws.Range("B22").Value = "Categorical"
$ws.Range("C22").Value = "Row Array"
$ws.Range("A23").Value = "order"
$ws.Range("B23").Value = "Categorical"
$ws.Range("C23").Value = "Cast Order"
$ws.Range("A24").Value = "director_ids"
$ws.Range("B24").Value = "Categorical"
$ws.Range("A25").Value = "IMDB Domestic Revenue"
$ws.Range("B25").Value = "Numerical"
$ws.Range("A26").Value = "international_revenue"
$ws.Range("B26").Value = "Numerical"
$ws.Range("A27").Value = "worldwide_revenue"
$ws.Range("B27").Value = "Numerical"
$ws.Range("A28").Value = "Release Year"
$ws.Range("B28").Value = "Categorical"
$ws.Range("A29").Value = "Merged Budget"
$ws.Range("B29").Value = "Categorical"
$ws.Range("C29").Value = "See Merge Methodology"
$ws.Range("A30").Value = "Merged Revenue"
$ws.Range("B30").Value = "Categorical"
$ws.Range("C30").Value = "See Merge Methodology"
$ws.Range("A31").Value = "Raw Star Scores"
$ws.Range("B31").Value = "Numerical"
$ws.Range("C31").Value = "Array"
$ws.Range("A32").Value = "Unweighted Star Score"
$ws.Range("B32").Value = "Numerical"
$ws.Range("C32").Value = "Uniform Sum"
$ws.Range("A33").Value = "Simple Weight Star Score"
$ws.Range("B33").Value = "Numerical"
$ws.Range("C33").Value = "Linear Weighted Sum"
$ws.Range("A34").Value = "Log Weight Star Score"
$ws.Range("B34").Value = "Numerical"
$ws.Range("C34").Value = "Log Weighted Sum"
$ws.Range("A35").Value = "Exponential Weight Star Score"
$ws.Range("B35").Value = "Numerical"
$ws.Range("C35").Value = "Exponential Weighted Sum"
$ws.Range("A36").Value = "director_scores"
$ws.Range("B36").Value = "Numerical"
$ws.Range("C36").Value = "Uniform Sum"
$ws.Range("A37").Value = "production_company_scores"
$ws.Range("B37").Value = "Numerical"
$ws.Range("C37").Value = "Uniform Sum"

# --- Update the view / selection to match ---
$ws.Range("B1:B37").Select()
